# 2A- Updated Script with order entry and Resouse sheet to
# Populates Sheet2 ("Resource") with order-entry data and makes it the
# active sheet/tab (moving tabSelected off Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- Header row -------------------------------------------------------
$headers = @("Customer ", "Order type", "PO", "Tag", "Delivery Term", `
  "Delivery Method", "Priority", "Item Line 1", "Order Qty", "Item Line 2", `
  "Order Qty", "Item Line 3", "Order qty", "Item Line 4", "Order Qty", `
  "Item Line 5", "Order Qty", "Blanket Agreement", "Promo", `
  "Transaction Reason", "Order count")

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# ---- Data rows ----------------------------------------------------------
# Columns used: A Customer, C Order type, D PO, E Tag, F Delivery Term,
# H Priority, I Order Qty, U Order count
$rows = @(
  @("US00025035", "04B WM Test", "Stock Balls", "PPD", "FXG",    "T2025S", 24,  1, $false),
  @("US00025035", "04B WM Test", "Stock Balls", "COL", "FXG",    "T2025S", 6,   1, $false),
  @("US00025025", "04B WM Test", "Stock Balls", "D75", "F02",    "T2025S", 4,   1, $false),
  @("US00002914", "04B WM Test", "Stock Balls", "PPD", "FXG",    "T2025S", 24,  1, $false),
  @("US00025687", "04B WM Test", "Stock Balls", "PPD", "F02",    "T2025S", 4,   1, $false),
  @("US00000215", "04B WM Test", "Stock Balls", "PPD", "F00",    "T2025S", 4,   1, $true),
  @("US00000215", "04B WM Test", "Stock Balls", "PPD", "F02",    "T2025S", 4,   1, $true),
  @("US00032279", "04B WM Test", "Stock Balls", "PPD", "F04",    "T2025S", 4,   1, $true),
  @("US00025033", "04B WM Test", "Stock Balls", "PPD", "FFE",    "T2025S", 504, 1, $true),
  @("US00025282", "04B WM Test", "Stock Balls", "PPD", "NMF",    "T2025S", 504, 1, $true)
)

$r = 2
foreach ($row in $rows) {
  $ws.Range("A$r").Value = $row[0]
  $ws.Range("C$r").Value = $row[1]
  $ws.Range("D$r").Value = $row[2]
  $ws.Range("E$r").Value = $row[3]

  $fCell = $ws.Range("F$r")
  if ($row[8]) {
    $fCell.NumberFormat = "0.00"
  }
  $fCell.Value = $row[4]

  $ws.Range("H$r").Value = $row[5]
  $ws.Range("I$r").Value = $row[6]
  $ws.Range("U$r").Value = $row[7]

  $r++
}

# ---- Column widths (best-fit sizing from the source workbook) ----------
$ws.Columns.Item(1).ColumnWidth = 11.28515625
$ws.Columns.Item(2).ColumnWidth = 10.5703125
$ws.Columns.Item(3).ColumnWidth = 12.28515625
$ws.Columns.Item(4).ColumnWidth = 10.28515625
$ws.Columns.Item(5).ColumnWidth = 13.5703125
$ws.Columns.Item(6).ColumnWidth = 16
$ws.Columns.Item(7).ColumnWidth = 7.5703125

# ---- Selection / activation --------------------------------------------
# Sheet2 becomes the active/selected tab (Sheet1 loses tabSelected).
[void]$ws.Range("F13").Select()
